$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header columns for Molar Volume in row 11 (the header row for the
# per-step results table), right after the existing "Fug. err (MPa)" column (Q11).
$ws.Range("R11").Value = "Molar Vol (L/mol)"
$ws.Range("S11").Value = "M. Vol. err (L/mol)"
